$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "'69.356.34"
$ws.Range("E2").Value = "  +0.63%  "

# Row 3
$ws.Range("D3").Formula = "'3.801.29"
$ws.Range("E3").Value = "  +1.34%  "

# Row 4
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Formula = "'604.32"
$ws.Range("E5").Value = "  +0.21%  "

# Row 6
$ws.Range("D6").Formula = "'165.31"
$ws.Range("E6").Value = "  -2.34%  "

# Row 7
$ws.Range("D7").Formula = "'3.799.71"
$ws.Range("E7").Value = "  +1.31%  "

# Row 9
$ws.Range("D9").Formula = "'0.540"
$ws.Range("E9").Value = "  +0.95%  "

# Row 10
$ws.Range("D10").Formula = "'0.172"
$ws.Range("E10").Value = "  +3.67%  "

# Row 11
$ws.Range("D11").Formula = "'6.34"
$ws.Range("E11").Value = "  -0.20%  "

# Row 12
$ws.Range("D12").Formula = "'0.463"
$ws.Range("E12").Value = "  -0.19%  "

# Row 13
$ws.Range("D13").Formula = "'37.47"
$ws.Range("E13").Value = "  -2.11%  "

# Row 14
$ws.Range("D14").Formula = "'0.0000248"
$ws.Range("E14").Value = "  -0.28%  "

# Row 15
$ws.Range("D15").Formula = "'4.439.98"
$ws.Range("E15").Value = "  +1.42%  "

# Row 16
$ws.Range("D16").Formula = "'3.795.24"
$ws.Range("E16").Value = "  +1.47%  "

# Row 17
$ws.Range("D17").Formula = "'69.499.46"

# Row 18
$ws.Range("D18").Formula = "'7.49"
$ws.Range("E18").Value = "  +2.58%  "

# Row 19
$ws.Range("D19").Formula = "'17.56"
$ws.Range("E19").Value = "  +2.94%  "

# Row 20
$ws.Range("E20").Value = "  -0.33%  "

# Row 21
$ws.Range("D21").Formula = "'11.20"
$ws.Range("E21").Value = "  +4.12%  "

# Row 22
$ws.Range("D22").Formula = "'494.18"
$ws.Range("E22").Value = "  -0.69%  "

# Row 23
$ws.Range("D23").Formula = "'0.726"
$ws.Range("E23").Value = "  -0.34%  "

# Row 24
$ws.Range("E24").Value = "  -1.74%  "

# Row 25
$ws.Range("D25").Formula = "'84.90"
$ws.Range("E25").Value = "  -0.68%  "

# Row 26
$ws.Range("E26").Value = "  -2.26%  "

# Row 27
$ws.Range("D27").Formula = "'12.32"
$ws.Range("E27").Value = "  -0.41%  "

# Row 28
$ws.Range("E28").Value = "  -1.83%  "

# Row 29
$ws.Range("E29").Value = "  +0.13%  "

# Row 30
$ws.Range("E30").Value = "  +0.24%  "

# Row 31
$ws.Range("E31").Value = "  +2.68%  "

# Row 32
$ws.Range("E32").Value = "  -4.49%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Formula = "'32.10"
$ws.Range("E33").Value = "  +0.63%  "

# Row 34
$ws.Range("B34").Value = "WrappedeETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D34").Formula = "'3.944.67"
$ws.Range("E34").Value = "  +1.28%  "

# Row 35
$ws.Range("D35").Formula = "'3.751.16"
$ws.Range("E35").Value = "  +1.81%  "

# Row 36
$ws.Range("E36").Value = "  -0.99%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Formula = "'0.141"
$ws.Range("E37").Value = "  +6.03%  "

# Row 38
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").Formula = "'1.02"
$ws.Range("E38").Value = "  -0.16%  "

# Row 39
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Formula = "'5.97"
$ws.Range("E39").Value = "  +1.77%  "

# Row 40
$ws.Range("D40").Formula = "'1.00"
$ws.Range("E40").Value = "  +0.00%  "

# Row 41
$ws.Range("E41").Value = "  +0.23%  "

# Row 42
$ws.Range("D42").Formula = "'3.06"
$ws.Range("E42").Value = "  +3.14%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Formula = "'1.99"
$ws.Range("E43").Value = "  +1.00%  "

# Row 44
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Formula = "'426.56"
$ws.Range("E44").Value = "  -2.69%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Formula = "'48.48"
$ws.Range("E45").Value = "  -0.87%  "

# Row 46
$ws.Range("D46").Formula = "'8.45"
$ws.Range("E46").Value = "  -0.44%  "

# Row 47
$ws.Range("E47").Value = "  -0.01%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Formula = "'142.46"
$ws.Range("E48").Value = "  +0.56%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Formula = "'2.819.94"
$ws.Range("E49").Value = "  +1.11%  "

# Row 50
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Formula = "'39.91"
$ws.Range("E50").Value = "  -1.60%  "

# Row 51
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Formula = "'1.29"
$ws.Range("E51").Value = "  +5.68%  "
